# Update the "ランサーズ" sheet with the new scrape snapshot taken at
# 2026-02-17 06:56:31. The previous snapshot had 18 listings (rows 2-19);
# the new one has 11 listings (rows 2-12), so the sheet shrinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2026-02-17 06:56:31"

# New listing data: Title, Category, Price, Deadline, URL, Score, Skills
$data = @(
    @("大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5423720", 385, "🔥AI,Ai ◆効率化"),
    @("建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5434128", 368, "🔥AI,Ai ◆開発"),
    @("企業のMicrosoft Copilot導入・活用支援AIコンサルタント募集(研修講師・メンター)", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5434363", 348, "🔥AI,Ai ◆コンサル"),
    @("【急募】製造業向け「製造副産物」の状態(硬度)判定AIのフィジビリティ検証(画像認識/動画解析)", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5439158", 303, "🔥AI,Ai"),
    @("※急募:Flutterによる業務アプリの開発(+next.js)", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5493471", 225, "🔥Next.js ◆開発 ◇アプリ"),
    @("※急募:Next.jsによる業務アプリの開発(+Flutter)", "システム開発", "300,000 円 ~ 500,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5493475", 225, "🔥Next.js ◆開発 ◇アプリ"),
    @("初回 【買い切り20万円】Shopeeチャット管理・返信Webツール開発(複数国対応)", "システム開発", "100,000 円 ~ 200,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5493016", 163, "◆ツール,開発 ◇管理"),
    @("【Unity/XRエンジニア募集】製造業DX支援!既存システムと連携するXRアプリ開発", "システム開発", "200,000 円 ~ 300,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5454210", 108, "◆開発 ◇アプリ"),
    @("【エクセル】教育機関向け教材販売・学習管理システムの構築(DB型設計・マトリックス集計)", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5493275", 48, "◇管理"),
    @("【急募】よもぎ蒸しサロンのWebサイトエラー解決依頼", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5493140", 33, "◇サイト"),
    @("【設計済み!作業時間~10時間】Stripe(銀行振込)を用いた月額課金システムの構築", "システム開発", "20,000 円 ~ 50,000 円 / 固定", "期限情報なし", "https://www.lancers.jp/work/detail/5493449", 28, "")
)

$oldLastRow = 19
$newLastRow = 1 + $data.Count   # 12

# Drop every existing hyperlink in the sheet up front (they'll be rebuilt
# below) so stale relationship entries don't linger once rows move/shrink.
$ws.Range("A1").Hyperlinks.Delete()

# Remove the rows that no longer exist in the new snapshot.
if ($oldLastRow -gt $newLastRow) {
    $clearRange = $ws.Range($ws.Cells.Item($newLastRow + 1, 1), $ws.Cells.Item($oldLastRow, 8))
    $clearRange.Clear()
}

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $rec = $data[$i]

    $ws.Cells.Item($row, 1).Value = $timestamp
    $ws.Cells.Item($row, 2).Value = $rec[0]
    $ws.Cells.Item($row, 3).Value = $rec[1]
    $ws.Cells.Item($row, 4).Value = $rec[2]
    $ws.Cells.Item($row, 5).Value = $rec[3]

    # URL column: real hyperlink + visible text, Hyperlink style applied
    # automatically by Hyperlinks.Add.
    $urlCell = $ws.Cells.Item($row, 6)
    $ws.Hyperlinks.Add($urlCell, $rec[4], "", "", $rec[4]) | Out-Null

    $ws.Cells.Item($row, 7).Value = $rec[5]
    $ws.Cells.Item($row, 8).Value = $rec[6]
}

# Column width adjustments
$ws.Columns.Item(2).ColumnWidth = 51
$ws.Columns.Item(4).ColumnWidth = 28
$ws.Columns.Item(8).ColumnWidth = 19
